$wb = $excel.ActiveWorkbook

# --- Update status text everywhere it appears ("In Translation" -> "Ready for handoff") ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("In Translation", "Ready for handoff") | Out-Null
}

# --- Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-04 18:43:51"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-04 18:43:46"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-04 18:43:51"

# --- Widen the status columns to fit the new "Ready for handoff" text ---
$overview.Columns.Item(5).ColumnWidth = 16.28
$overview.Columns.Item(6).ColumnWidth = 16.28
$zhcn.Columns.Item(3).ColumnWidth = 16.28
$dede.Columns.Item(3).ColumnWidth = 16.28
